$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: column letter, row number, new value.
# Column D values are prefixed with a leading apostrophe so Excel keeps
# them stored as text (matching the source inlineStr cells) instead of
# auto-converting numeric-looking strings (e.g. "0.999") into numbers.
$updates = @(
    @("D", 2, "'64.799.78"),
    @("E", 2, "  -0.63%  "),
    @("D", 3, "'3.143.45"),
    @("E", 3, "  -1.12%  "),
    @("D", 4, "'0.999"),
    @("E", 4, "  -0.20%  "),
    @("D", 5, "'581.07"),
    @("E", 5, "  +1.12%  "),
    @("D", 6, "'147.47"),
    @("E", 6, "  -2.53%  "),
    @("E", 7, "  +0.07%  "),
    @("D", 8, "'3.142.16"),
    @("E", 8, "  -1.10%  "),
    @("D", 9, "'0.525"),
    @("E", 9, "  -0.71%  "),
    @("D", 10, "'0.158"),
    @("E", 10, "  -3.65%  "),
    @("D", 11, "'6.17"),
    @("E", 11, "  -1.35%  "),
    @("D", 12, "'0.497"),
    @("E", 12, "  -2.51%  "),
    @("E", 13, "  -2.50%  "),
    @("D", 14, "'37.12"),
    @("E", 14, "  -3.20%  "),
    @("D", 15, "'3.652.75"),
    @("E", 15, "  -1.32%  "),
    @("D", 16, "'64.804.51"),
    @("E", 16, "  -0.79%  "),
    @("B", 17, "Polkadot"),
    @("C", 17, "https://coinranking.com/coin/25W7FG7om+polkadot-dot"),
    @("D", 17, "'7.13"),
    @("E", 17, "  -1.53%  "),
    @("B", 18, "WrappedEther"),
    @("C", 18, "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"),
    @("D", 18, "'3.137.89"),
    @("E", 18, "  -1.11%  "),
    @("E", 19, "  -0.54%  "),
    @("D", 20, "'499.54"),
    @("E", 20, "  -2.87%  "),
    @("D", 21, "'15.33"),
    @("E", 21, "  +2.02%  "),
    @("D", 22, "'0.712"),
    @("E", 22, "  -3.89%  "),
    @("D", 23, "'15.00"),
    @("E", 23, "  -7.39%  "),
    @("D", 24, "'7.76"),
    @("E", 24, "  -1.76%  "),
    @("D", 25, "'84.20"),
    @("E", 25, "  -1.03%  "),
    @("D", 27, "'9.13"),
    @("E", 27, "  -1.07%  "),
    @("E", 28, "  -0.21%  "),
    @("D", 29, "'2.19"),
    @("E", 29, "  -1.23%  "),
    @("D", 30, "'2.82"),
    @("E", 30, "  +1.10%  "),
    @("D", 31, "'27.54"),
    @("E", 31, "  -2.00%  "),
    @("E", 32, "  -0.78%  "),
    @("E", 33, "  +0.02%  "),
    @("E", 34, "  +1.04%  "),
    @("D", 35, "'6.45"),
    @("E", 35, "  -3.44%  "),
    @("D", 36, "'54.86"),
    @("E", 36, "  -1.81%  "),
    @("D", 37, "'0.0894"),
    @("E", 37, "  +1.68%  "),
    @("D", 38, "'470.11"),
    @("E", 38, "  -1.94%  "),
    @("E", 39, "  -1.57%  "),
    @("D", 40, "'2.92"),
    @("E", 40, "  -7.58%  "),
    @("D", 41, "'8.73"),
    @("E", 41, "  +0.45%  "),
    @("D", 42, "'2.977.69"),
    @("E", 42, "  -4.76%  "),
    @("E", 43, "  -4.25%  "),
    @("D", 44, "'2.43"),
    @("E", 44, "  -4.53%  "),
    @("D", 45, "'0.282"),
    @("E", 45, "  -3.64%  "),
    @("D", 46, "'28.24"),
    @("E", 46, "  -4.05%  "),
    @("D", 47, "'0.0₃0600"),
    @("E", 47, "  +1.64%  "),
    @("E", 49, "  -1.84%  "),
    @("D", 50, "'2.23"),
    @("E", 50, "  -5.06%  "),
    @("D", 51, "'118.95"),
    @("E", 51, "  -4.67%  ")
)

foreach ($u in $updates) {
    $col = $u[0]
    $row = $u[1]
    $val = $u[2]
    $addr = "$col$row"
    $ws.Range($addr).Value2 = $val
}
